$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), extending the header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, thin border) from the existing
# header cell H1 onto the two new header cells so they reuse the same style
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New numeric data cells for rows 2 and 3
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
